$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "doOutsCrop" row (old row 9); later rows shift up by one.
$ws.Rows.Item(9).Delete()

# Replace the "ru" language column header/values with the new German ("de") column.
$ws.Range("D1").Value = "de"
$ws.Range("D2").Value = "Wähle Schnittstil"
$ws.Range("D3").Value = "Wähle Schnittstil"
$ws.Range("D4").Value = "Arbeitsfläche aufteilen"
$ws.Range("D5").Value = "Erstelle Schnittmaske"
$ws.Range("D6").Value = "Abbrechen"
$ws.Range("D7").Value = "Hintergrund auf Ebene"
$ws.Range("D8").Value = "Hintergrund füllen"
$ws.Range("D9").Value = " - Rastern"
$ws.Range("D10").Value = " - Skalieren"
$ws.Range("D11").Value = " - Einblenden"
$ws.Range("D12").Value = " - Beschneiden"
$ws.Range("D13").Value = "Golden Crop by SzopeN"
$ws.Range("D14").Value = "Schnittmaske"
$ws.Range("D15").Value = "Trennungsregeln"
$ws.Range("D16").Value = "Linien auf %1%%"
$ws.Range("D17").Value = "Goldene Diagonale aufwärts"
$ws.Range("D18").Value = "Goldene Diagonale abwärts"
$ws.Range("D19").Value = "Öffne das Dokument, in dem das Script ablaufen soll."
$ws.Range("D20").Value = "Erweiterung der Arbeitsfläche zeigen"
$ws.Range("D21").Value = "Was mache ich mit der Arbeitsfläche?"
$ws.Range("D22").Value = "Erweiterte Arbeitsfläche"
$ws.Range("D23").Value = "Schnitt ohne Erweiterung"
$ws.Range("D24").Value = "Zurück zum Schneiden"
$ws.Range("D25").Value = "Kompositionsmethode"
$ws.Range("D26").Value = "Auswahl der Kompositionslinien"
$ws.Range("D27").Value = "Goldene Regel"
$ws.Range("D28").Value = "Drittel-Regel"
$ws.Range("D29").Value = "Goldene Spirale unten links"
$ws.Range("D30").Value = "Goldene Spirale oben links"
$ws.Range("D31").Value = "Goldene Spirale oben rechts"
$ws.Range("D32").Value = "Goldene Spirale unten rechts"
$ws.Range("D33").Value = "Alles Auswählen"
$ws.Range("D34").Value = "Nichts Auswählen"
$ws.Range("D35").Value = "OK"
$ws.Range("D36").Value = "Alle Goldenen Spiralen"
$ws.Range("D37").Value = "Grundregeln"

# Widen column D to fit the German text, matching the authored layout.
$ws.Columns.Item(4).ColumnWidth = 72.16666666666667

# Reset the view: zoom out and select the top-left cell (matches the commit's saved view state).
$excel.ActiveWindow.Zoom = 85
$ws.Range("A1").Select()
